$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 168201
$ws.Range("C4").Value = 159076
$ws.Range("C5").Value = 9125
$ws.Range("C8").Value = 65.45999999999999
